$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing row 10 (WS-2412-008 / John_Doe / 12/18/2024 / Approved) as a
# formatting template for the new row 11 so the new cells inherit the same
# cell style (s="1") instead of Excel auto-detecting the date text and minting
# a brand-new number-format style.
$ws.Range("A10:E10").Copy($ws.Range("A11:E11"))

$ws.Range("A11").Value = "WS-2412-009"
$ws.Range("B11").Value = "Jane"
$ws.Range("D11").Formula = '=HYPERLINK("C:/Users/admin/git-directories/Data_Management/Inventory_Management/Withdrawals/WS-2412-009_Jane.pdf", "WS-2412-009_Jane.pdf")'
# C11 ("12/18/2024" text) and E11 ("Approved") already carry over correctly
# from the template row, matching the source data for the new withdrawal entry.

$excel.CutCopyMode = 0
